# Applies the "Updated symbol list" edit described in the commit diff.
# Columns: A=index, B=Coin, C=Link, D=Price, E=RankCoinSymbol, F=Date, G=Hour
# Rows 18-24 show coins shifting up by one position (a cyclic re-rank),
# with row 24 wrapping to the coin that used to be in row 18 ("One").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All Price (column D) cells are stored as text in the workbook (e.g. "0.1330"
# keeps its trailing zero). Writing a leading apostrophe forces Excel to keep the
# assigned text verbatim instead of re-interpreting it as a number, and resetting
# the Style back to "Normal" avoids Excel applying an automatic "Text" number
# format/style to the cell as a side effect of the apostrophe.
function Set-TextPrice($cellRef, $text) {
    $ws.Range($cellRef).Value = "'" + $text
    $ws.Range($cellRef).Style = "Normal"
}

# Row 2
Set-TextPrice "D2" "245.58"
# Row 3
Set-TextPrice "D3" "23.79"
# Row 4
Set-TextPrice "D4" "5.344"
# Row 5
Set-TextPrice "D5" "0.05838"
# Row 6
Set-TextPrice "D6" "6.479"
# Row 8
Set-TextPrice "D8" "0.8117"
# Row 9
Set-TextPrice "D9" "0.9223"
# Row 11
Set-TextPrice "D11" "0.07366"
# Row 12
Set-TextPrice "D12" "0.03096"
# Row 13
Set-TextPrice "D13" "0.03055"
# Row 14
Set-TextPrice "D14" "0.09374"
# Row 15
Set-TextPrice "D15" "3.858"
# Row 16
Set-TextPrice "D16" "0.001548"
# Row 17
Set-TextPrice "D17" "0.04687"
# Row 18
$ws.Range("B18").Value = "TigerCash"
$ws.Range("C18").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextPrice "D18" "0.005989"
$ws.Range("E18").Value = "17TigerCashTCH"
# Row 19
$ws.Range("B19").Value = "BitKan"
$ws.Range("C19").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
Set-TextPrice "D19" "0.001251"
$ws.Range("E19").Value = "18BitKanKAN"
# Row 20
$ws.Range("B20").Value = "HotbitToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
Set-TextPrice "D20" "0.004687"
$ws.Range("E20").Value = "19HotbitTokenHTB"
# Row 21
$ws.Range("B21").Value = "NitroEx"
$ws.Range("C21").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
Set-TextPrice "D21" "0.00008819"
$ws.Range("E21").Value = "20NitroExNTX"
# Row 22
$ws.Range("B22").Value = "LEO"
$ws.Range("C22").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextPrice "D22" "3.592"
$ws.Range("E22").Value = "21LEOLEO"
# Row 23
$ws.Range("B23").Value = "BTSEToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextPrice "D23" "2.158"
$ws.Range("E23").Value = "22BTSETokenBTSE"
# Row 24
$ws.Range("B24").Value = "One"
$ws.Range("C24").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextPrice "D24" "0.01085"
$ws.Range("E24").Value = "23OneONEBestin24h"
# Row 26
Set-TextPrice "D26" "0.1330"
# Row 40
Set-TextPrice "D40" "0.03848"
# Row 41
Set-TextPrice "D41" "0.006414"
# Row 42
Set-TextPrice "D42" "0.1066"
# Row 43
Set-TextPrice "D43" "0.002914"
# Row 44
Set-TextPrice "D44" "0.008505"
# Row 45
Set-TextPrice "D45" "0.00005255"
# Row 47
Set-TextPrice "D47" "0.6534"
# Row 48
Set-TextPrice "D48" "0.001863"
# Row 50
Set-TextPrice "D50" "0.0002004"
